$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '70.468.43'
$ws.Range("E2").Value = '  +0.94%  '

# Row 3
$ws.Range("D3").Value = '3.520.90'
$ws.Range("E3").Value = '  +0.19%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.81'
$ws.Range("E5").Value = '  +0.51%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.49'
$ws.Range("E6").Value = '  +1.98%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.613'
$ws.Range("E7").Value = '  -0.36%  '

# Row 8
$ws.Range("D8").Value = '3.516.32'
$ws.Range("E8").Value = '  +0.18%  '

# Row 10
$ws.Range("E10").Value = '  -1.24%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.25'
$ws.Range("E11").Value = '  +8.65%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.588'
$ws.Range("E12").Value = '  +0.99%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.36'
$ws.Range("E13").Value = '  -1.77%  '

# Row 14
$ws.Range("E14").Value = '  -0.19%  '

# Row 15
$ws.Range("D15").Value = '4.092.31'
$ws.Range("E15").Value = '  +0.19%  '

# Row 16
$ws.Range("E16").Value = '  -0.60%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '613.26'

# Row 18
$ws.Range("D18").Value = '3.524.24'
$ws.Range("E18").Value = '  +0.37%  '

# Row 19
$ws.Range("D19").Value = '70.509.20'
$ws.Range("E19").Value = '  +0.99%  '

# Row 20
$ws.Range("E20").Value = '  +0.96%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.52'
$ws.Range("E21").Value = '  +1.36%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.882'
$ws.Range("E22").Value = '  -0.19%  '

# Row 23
$ws.Range("E23").Value = '  -9.16%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '99.42'
$ws.Range("E24").Value = '  +3.76%  '

# Row 25
$ws.Range("E25").Value = '  -0.55%  '

# Row 26
$ws.Range("E26").Value = '  -2.76%  '

# Row 27
$ws.Range("E27").Value = '  -0.09%  '

# Row 28
$ws.Range("E28").Value = '  -0.79%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.33'
$ws.Range("E29").Value = '  +3.56%  '

# Row 30
$ws.Range("E30").Value = '  -1.62%  '

# Row 31
$ws.Range("E31").Value = '  -4.05%  '

# Row 32
$ws.Range("E32").Value = '  -2.76%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '645.71'
$ws.Range("E33").Value = '  +13.66%  '

# Row 34
$ws.Range("B34").Value = 'Mantle'
$ws.Range("C34").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.28'
$ws.Range("E34").Value = '  -4.10%  '

# Row 35
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.85'
$ws.Range("E35").Value = '  -1.58%  '

# Row 36
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0998'
$ws.Range("E36").Value = '  -1.27%  '

# Row 37
$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.56'
$ws.Range("E37").Value = '  +2.53%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.78'
$ws.Range("E38").Value = '  +0.22%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0478'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '56.83'
$ws.Range("E40").Value = '  -0.38%  '

# Row 41
$ws.Range("E41").Value = '  +0.10%  '

# Row 42
$ws.Range("E42").Value = '  +1.36%  '

# Row 43
$ws.Range("E43").Value = '  +6.69%  '

# Row 44
$ws.Range("D44").Value = '3.367.99'
$ws.Range("E44").Value = '  +1.05%  '

# Row 45
$ws.Range("E45").Value = '  -4.68%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '32.26'
$ws.Range("E46").Value = '  -2.49%  '

# Row 47
$ws.Range("E47").Value = '  -3.25%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.56'
$ws.Range("E48").Value = '  -2.49%  '

# Row 49
$ws.Range("E49").Value = '  +1.01%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.43'
$ws.Range("E50").Value = '  -1.64%  '

# Row 51
$ws.Range("E51").Value = '  -0.01%  '
